$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.5333333333333333, 0.6153846153846154, 0.5714285714285715, 13)
    3  = @(0.4444444444444444, 0.3636363636363636, 0.4, 11)
    4  = @(0.5, 0.5, 0.5, 0.5)
    5  = @(0.4888888888888889, 0.4895104895104895, 0.4857142857142858, 24)
    6  = @(0.4925925925925926, 0.5, 0.4928571428571429, 24)
    7  = @(0.5625, 0.6923076923076923, 0.6206896551724138, 13)
    8  = @(0.5, 0.3636363636363636, 0.4210526315789474, 11)
    9  = @(0.5416666666666666, 0.5416666666666666, 0.5416666666666666, 0.5416666666666666)
    10 = @(0.53125, 0.5279720279720279, 0.5208711433756806, 24)
    11 = @(0.5338541666666666, 0.5416666666666666, 0.529189352692075, 24)
    12 = @(0.5454545454545454, 0.4615384615384616, 0.4999999999999999, 13)
    13 = @(0.4615384615384616, 0.5454545454545454, 0.4999999999999999, 11)
    14 = @(0.5, 0.5, 0.5, 0.5)
    15 = @(0.5034965034965035, 0.5034965034965035, 0.4999999999999999, 24)
    16 = @(0.506993006993007, 0.5, 0.4999999999999998, 24)
    17 = @(0.6153846153846154, 0.6153846153846154, 0.6153846153846154, 13)
    18 = @(0.5454545454545454, 0.5454545454545454, 0.5454545454545454, 11)
    19 = @(0.5833333333333334, 0.5833333333333334, 0.5833333333333334, 0.5833333333333334)
    20 = @(0.5804195804195804, 0.5804195804195804, 0.5804195804195804, 24)
    21 = @(0.5833333333333334, 0.5833333333333334, 0.5833333333333334, 24)
    22 = @(0.5, 0.5384615384615384, 0.5185185185185186, 13)
    23 = @(0.4, 0.3636363636363636, 0.380952380952381, 11)
    24 = @(0.4583333333333333, 0.4583333333333333, 0.4583333333333333, 0.4583333333333333)
    25 = @(0.45, 0.451048951048951, 0.4497354497354498, 24)
    26 = @(0.4541666666666667, 0.4583333333333333, 0.4554673721340389, 24)
}

$cols = @("B", "C", "D", "E")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
